# Reproduce the workbook edit described by the fixture diff:
#   - Tests!B5 changes from the number 0 to the text "string"
#   - Tests!C6 changes from the text "'" to the text "string"
#   - Tests row 5 / row 6 heights change to 14.9 / 15
#   - The active/selected sheet moves from "Main root" to "Tests"
#   - "Main root" keeps cursor position A2 on-screen, but the saved
#     selection on its frozen bottom-right pane moves to C4
#   - "Tests" sheet retains its existing bottom-right selection (C6)

$wb = $excel.ActiveWorkbook

$mainRoot = $wb.Worksheets.Item(1)
$tests = $wb.Worksheets.Item(2)

# --- data edits on the "Tests" sheet -------------------------------------
$tests.Range("B5").Value = "string"
$tests.Range("C6").Value = "string"

# --- row height tweaks on the "Tests" sheet ------------------------------
$tests.Rows.Item(5).RowHeight = 14.9
$tests.Rows.Item(6).RowHeight = 15

# --- selection / active-sheet changes ------------------------------------
# Move the cursor on "Main root" (updates its saved pane selection to C4)
$mainRoot.Range("C4").Select()

# Make "Tests" the active/selected sheet (updates workbook activeTab and
# each sheet's tabSelected flag)
$tests.Activate()
